$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.332.88"
$ws.Range("E2").Value = '  +2.53%  '

$ws.Range("D3").Value = "'2.427.89"
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("D5").Value = "'318.83"
$ws.Range("E5").Value = '  +3.46%  '

$ws.Range("D6").Value = "'102.76"
$ws.Range("E6").Value = '  +1.99%  '

$ws.Range("D7").Value = "'0.517"
$ws.Range("E7").Value = '  +0.45%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = '  +5.38%  '

$ws.Range("D10").Value = "'35.63"
$ws.Range("E10").Value = '  +0.53%  '

$ws.Range("D11").Value = "'0.0801"
$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("E12").Value = '  -1.99%  '

$ws.Range("D13").Value = "'18.24"
$ws.Range("E13").Value = '  -3.05%  '

$ws.Range("D14").Value = "'7.09"
$ws.Range("E14").Value = '  +2.18%  '

$ws.Range("D15").Value = "'2.806.94"
$ws.Range("E15").Value = '  +0.02%  '

$ws.Range("D16").Value = "'2.418.66"
$ws.Range("E16").Value = '  +0.20%  '

$ws.Range("D17").Value = "'0.847"
$ws.Range("E17").Value = '  +1.24%  '

$ws.Range("D18").Value = "'45.204.68"
$ws.Range("E18").Value = '  +2.31%  '

$ws.Range("D19").Value = "'12.23"
$ws.Range("E19").Value = '  -0.87%  '

$ws.Range("D20").Value = "'6.35"
$ws.Range("E20").Value = '  -1.45%  '

$ws.Range("E21").Value = '  +1.63%  '

$ws.Range("D22").Value = "'68.91"
$ws.Range("E22").Value = '  +0.47%  '

$ws.Range("D23").Value = "'244.86"
$ws.Range("E23").Value = '  +1.66%  '

$ws.Range("D24").Value = "'2.28"
$ws.Range("E24").Value = '  -1.01%  '

$ws.Range("E25").Value = '  +0.17%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").Value = "'25.81"
$ws.Range("E27").Value = '  +2.08%  '

$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = '  -2.57%  '

$ws.Range("D29").Value = "'9.62"
$ws.Range("E29").Value = '  +0.25%  '

$ws.Range("D30").Value = "'49.48"
$ws.Range("E30").Value = '  +2.60%  '

$ws.Range("D31").Value = "'33.02"
$ws.Range("E31").Value = '  -0.85%  '

$ws.Range("D32").Value = "'0.125"
$ws.Range("E32").Value = '  +5.82%  '

$ws.Range("D33").Value = "'20.01"
$ws.Range("E33").Value = '  +7.23%  '

$ws.Range("E34").Value = '  +0.66%  '

$ws.Range("E35").Value = '  +0.25%  '

$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("E37").Value = '  -1.98%  '

$ws.Range("E38").Value = '  -0.91%  '

$ws.Range("D39").Value = "'2.88"
$ws.Range("E39").Value = '  -1.25%  '

$ws.Range("D40").Value = "'126.32"
$ws.Range("E40").Value = '  -2.85%  '

$ws.Range("E41").Value = '  +0.75%  '

$ws.Range("E42").Value = '  -3.94%  '

$ws.Range("D43").Value = "'20.73"
$ws.Range("E43").Value = '  -1.57%  '

$ws.Range("E44").Value = '  +0.64%  '

$ws.Range("D45").Value = "'1.937.30"
$ws.Range("E45").Value = '  -1.08%  '

$ws.Range("E46").Value = '  -2.62%  '

$ws.Range("E47").Value = '  +2.09%  '

$ws.Range("D48").Value = "'1.81"
$ws.Range("E48").Value = '  +9.32%  '

$ws.Range("E49").Value = '  -3.47%  '

$ws.Range("D50").Value = "'76.81"
$ws.Range("E50").Value = '  +4.47%  '

$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = "'4.79"
$ws.Range("E51").Value = '  +5.33%  '
